# ExcelDatosCuentasValidaDNI.xlsx - regression R33 data update
# - updates the DNI/document number and the trailing sequence number on
#   row 11 of Hoja1
# - moves the viewport/selection so N11 is the active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Data changes: row 11 (Smoke user) gets a new DNI and a new running number
$ws.Range("G11").Value = 24741865
$ws.Range("N11").Value = 307

# View changes: scroll the window left a column and move the selection to N11
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("N11").Select()
